$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 371.85715
$ws.Range("I9").Value = 405.625
$ws.Range("J9").Value = 263.8
$ws.Range("K9").Value = 405.625
$ws.Range("L9").Value = 263.8
$ws.Range("M9").Value = -236.625
$ws.Range("N9").Value = -601.8
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H38").Value = 4146.7856
$ws.Range("I38").Value = 1005.1818
$ws.Range("J38").Value = 15666
$ws.Range("K38").Value = 3015.5454
$ws.Range("L38").Value = 46998
$ws.Range("M38").Value = -2643.5454
$ws.Range("N38").Value = -47742
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H100").Value = 2200
$ws.Range("I100").Value = 1983.3334
$ws.Range("K100").Value = 1983.3334
$ws.Range("M100").Value = -1442.3334
$ws.Range("H101").Value = 3788697
$ws.Range("J101").Value = 1093.1666
$ws.Range("L101").Value = 3279.4998
$ws.Range("N101").Value = -6523.4998
$ws.Range("H111").Value = 5188
$ws.Range("I111").Value = 4032.8333
$ws.Range("K111").Value = 12098.4999
$ws.Range("M111").Value = -9031.499899999999
$ws.Range("H115").Value = 1951699.5
$ws.Range("I115").Value = 2059849.5
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 6179548.5
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -6177981.5
$ws.Range("N115").Value = -18134
$ws.Range("H116").Value = 11444.4375
$ws.Range("I116").Value = 12412.786
$ws.Range("K116").Value = 12412.786
$ws.Range("M116").Value = -8970.786
$ws.Range("H125").Value = 5529038.5
$ws.Range("I125").Value = 1359.25
$ws.Range("J125").Value = 11056718
$ws.Range("K125").Value = 12233.25
$ws.Range("L125").Value = 99510462
$ws.Range("M125").Value = -9773.25
$ws.Range("N125").Value = -99515382
$ws.Range("H135").Value = 1177.1714
$ws.Range("I135").Value = 804.5
$ws.Range("K135").Value = 7240.5
$ws.Range("M135").Value = -4705.5
$ws.Range("H137").Value = 2939.1428
$ws.Range("I137").Value = 2713.85
$ws.Range("J137").Value = 3330.9565
$ws.Range("K137").Value = 8141.549999999999
$ws.Range("L137").Value = 9992.869499999999
$ws.Range("M137").Value = -5591.549999999999
$ws.Range("N137").Value = -15092.8695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 6593.1685
$ws.Range("I32").Value = 5057.5684
$ws.Range("K32").Value = 5057.5684
$ws.Range("M32").Value = -4770.5684
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 3688.0667
$ws.Range("I74").Value = 3141.2307
$ws.Range("K74").Value = 3141.2307
$ws.Range("M74").Value = -2267.2307
$ws.Range("H77").Value = 3688.0667
$ws.Range("I77").Value = 3141.2307
$ws.Range("K77").Value = 15706.1535
$ws.Range("M77").Value = -11338.1535
$ws.Range("H122").Value = 6647.7
$ws.Range("I122").Value = 2615.25
$ws.Range("K122").Value = 7845.75
$ws.Range("M122").Value = -5395.75
$ws.Range("H132").Value = 2464.1277
$ws.Range("I132").Value = 2599.1707
$ws.Range("K132").Value = 7797.5121
$ws.Range("M132").Value = -5267.5121
$ws.Range("H139").Value = 42499.5
$ws.Range("J139").Value = 42499.5
$ws.Range("L139").Value = 42499.5
$ws.Range("N139").Value = -52779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1958
$ws.Range("J64").Value = 2102.111
$ws.Range("L64").Value = 2102.111
$ws.Range("N64").Value = -2552.111
$ws.Range("H67").Value = 1958
$ws.Range("J67").Value = 2102.111
$ws.Range("L67").Value = 2102.111
$ws.Range("N67").Value = -3662.111
$ws.Range("H107").Value = 1110.4546
$ws.Range("J107").Value = 1002.4
$ws.Range("L107").Value = 1002.4
$ws.Range("N107").Value = -4842.4
$ws.Range("H134").Value = 3305.476
$ws.Range("I134").Value = 2649.4
$ws.Range("J134").Value = 4945.6665
$ws.Range("K134").Value = 7948.200000000001
$ws.Range("L134").Value = 14836.9995
$ws.Range("M134").Value = -5413.200000000001
$ws.Range("N134").Value = -19906.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2653
$ws.Range("I16").Value = 1742.5714
$ws.Range("J16").Value = 4246.25
$ws.Range("K16").Value = 1742.5714
$ws.Range("L16").Value = 4246.25
$ws.Range("M16").Value = -1455.5714
$ws.Range("N16").Value = -4820.25
$ws.Range("H31").Value = 2149.5945
$ws.Range("I31").Value = 1783.1578
$ws.Range("J31").Value = 2536.389
$ws.Range("K31").Value = 1783.1578
$ws.Range("L31").Value = 2536.389
$ws.Range("M31").Value = -1488.1578
$ws.Range("N31").Value = -3126.389
$ws.Range("H34").Value = 2149.5945
$ws.Range("I34").Value = 1783.1578
$ws.Range("J34").Value = 2536.389
$ws.Range("K34").Value = 1783.1578
$ws.Range("L34").Value = 2536.389
$ws.Range("M34").Value = -1581.1578
$ws.Range("N34").Value = -2940.389
$ws.Range("H86").Value = 4702.353
$ws.Range("I86").Value = 3616.2222
$ws.Range("J86").Value = 5924.25
$ws.Range("K86").Value = 3616.2222
$ws.Range("L86").Value = 5924.25
$ws.Range("M86").Value = -2493.2222
$ws.Range("N86").Value = -8170.25
$ws.Range("H89").Value = 4702.353
$ws.Range("I89").Value = 3616.2222
$ws.Range("J89").Value = 5924.25
$ws.Range("K89").Value = 18081.111
$ws.Range("L89").Value = 29621.25
$ws.Range("M89").Value = -12465.111
$ws.Range("N89").Value = -40853.25
$ws.Range("H107").Value = 884.5
$ws.Range("I107").Value = 481.33334
$ws.Range("K107").Value = 481.33334
$ws.Range("M107").Value = 1438.66666
$ws.Range("H113").Value = 2653
$ws.Range("I113").Value = 1742.5714
$ws.Range("J113").Value = 4246.25
$ws.Range("K113").Value = 1742.5714
$ws.Range("L113").Value = 4246.25
$ws.Range("M113").Value = 427.4286
$ws.Range("N113").Value = -8586.25
$ws.Range("H132").Value = 1674.5333
$ws.Range("I132").Value = 1689.1316
$ws.Range("J132").Value = 1595.2858
$ws.Range("K132").Value = 5067.3948
$ws.Range("L132").Value = 4785.857400000001
$ws.Range("M132").Value = -2537.3948
$ws.Range("N132").Value = -9845.857400000001
$ws.Range("H134").Value = 866.05884
$ws.Range("I134").Value = 817.2286
$ws.Range("K134").Value = 2451.6858
$ws.Range("M134").Value = 83.3141999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 632.2105
$ws.Range("I12").Value = 828.4
$ws.Range("J12").Value = 562.1429000000001
$ws.Range("K12").Value = 2485.2
$ws.Range("L12").Value = 1686.4287
$ws.Range("M12").Value = -2312.2
$ws.Range("N12").Value = -2032.4287
$ws.Range("H38").Value = 202.125
$ws.Range("I38").Value = 173.73334
$ws.Range("J38").Value = 249.44444
$ws.Range("K38").Value = 521.20002
$ws.Range("L38").Value = 748.33332
$ws.Range("M38").Value = -174.20002
$ws.Range("N38").Value = -1442.33332
$ws.Range("H92").Value = 658.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9285.799999999999
$ws.Range("J70").Value = 8785
$ws.Range("L70").Value = 8785
$ws.Range("N70").Value = -9325
$ws.Range("H73").Value = 9285.799999999999
$ws.Range("J73").Value = 8785
$ws.Range("L73").Value = 8785
$ws.Range("N73").Value = -10657
$ws.Range("H97").Value = 1974.6863
$ws.Range("I97").Value = 1543.9706
$ws.Range("K97").Value = 1543.9706
$ws.Range("M97").Value = -1047.9706
$ws.Range("H122").Value = 2212.5715
$ws.Range("I122").Value = 2498.25
$ws.Range("K122").Value = 7494.75
$ws.Range("M122").Value = -5044.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7419.448
$ws.Range("I7").Value = 7574.2593
$ws.Range("K7").Value = 7574.2593
$ws.Range("M7").Value = -7462.2593
$ws.Range("H40").Value = 10512.4
$ws.Range("I40").Value = 12987.682
$ws.Range("J40").Value = 3705.375
$ws.Range("K40").Value = 12987.682
$ws.Range("L40").Value = 3705.375
$ws.Range("M40").Value = -12851.682
$ws.Range("N40").Value = -3977.375
$ws.Range("H126").Value = 7419.448
$ws.Range("I126").Value = 7574.2593
$ws.Range("K126").Value = 22722.7779
$ws.Range("M126").Value = -20252.7779

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17355.75
$ws.Range("J41").Value = 17355.75
$ws.Range("L41").Value = 17355.75
$ws.Range("N41").Value = -18135.75
$ws.Range("H107").Value = 937.4583
$ws.Range("I107").Value = 890.9
$ws.Range("K107").Value = 2672.7
$ws.Range("M107").Value = -752.6999999999998
$ws.Range("H136").Value = 3110.8677
$ws.Range("I136").Value = 2543.712
$ws.Range("J136").Value = 6828.8887
$ws.Range("K136").Value = 7631.136
$ws.Range("L136").Value = 20486.6661
$ws.Range("M136").Value = -5081.136
$ws.Range("N136").Value = -25586.6661
